$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column F ("Popis Cesty") before the existing "Stav Tachometra" column ---
$ws.Columns.Item(6).Insert()

# --- Header row (row 1) ---
$ws.Cells.Item(1, 6).Value = "Popis Cesty"
$ws.Cells.Item(1, 7).Value = "Stav Tachometra"
$ws.Cells.Item(1, 8).Value = "Km Jazda"

# Copy the header style (font/border/alignment) from an existing styled header
# cell onto the new F1 header cell so it matches the rest of the header row.
$ws.Cells.Item(1, 1).Copy()
$ws.Cells.Item(1, 6).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows (row 2 .. row 19) ---
$data = @(
    @("2025-11-03", "Vrbové",         "07:00", "Bratislava",      "08:00", "Obchodné rokovanie o IT",        125746, 92),
    @("2025-11-03", "Bratislava",     "15:30", "Vrbové",          "16:30", "Obchodné rokovanie o IT",        125838, 92),
    @("2025-11-04", "Vrbové",         "06:30", "Trnava",          "07:30", "Kontrola technického vybavenia", 125886, 48),
    @("2025-11-04", "Trnava",         "15:00", "Vrbové",          "16:00", "Kontrola technického vybavenia", 125934, 48),
    @("2025-11-05", "Vrbové",         "07:30", "Piešťany",        "08:30", "Servis IT infraštruktúry",       125945, 11),
    @("2025-11-05", "Piešťany",       "16:00", "Vrbové",          "17:00", "Servis IT infraštruktúry",       125956, 11),
    @("2025-11-06", "Vrbové",         "06:45", "Nitra",           "07:45", "Konzultácia vývoja softvéru",    126016, 60),
    @("2025-11-06", "Nitra",          "15:15", "Vrbové",          "16:15", "Konzultácia vývoja softvéru",    126076, 60),
    @("2025-11-07", "Vrbové",         "07:15", "Trenčín",         "08:15", "Implementácia cloud riešenia",   126132, 56),
    @("2025-11-07", "Trenčín",        "16:00", "Vrbové",          "17:00", "Implementácia cloud riešenia",   126188, 56),
    @("2025-11-10", "Vrbové",         "06:00", "Žilina",          "07:00", "Analýza bezpečnostných rizík",   126314, 126),
    @("2025-11-10", "Žilina",         "17:00", "Vrbové",          "18:00", "Analýza bezpečnostných rizík",   126440, 126),
    @("2025-11-11", "Vrbové",         "06:30", "Banská Bystrica", "07:30", "Obchodné rokovanie o IT",        126642, 202),
    @("2025-11-11", "Banská Bystrica","18:00", "Vrbové",          "19:00", "Obchodné rokovanie o IT",        126844, 202),
    @("2025-11-12", "Vrbové",         "07:00", "Martin",          "08:00", "Kontrola technického vybavenia", 126928, 84),
    @("2025-11-12", "Martin",         "16:00", "Vrbové",          "17:00", "Kontrola technického vybavenia", 127012, 84),
    @("2025-11-14", "Vrbové",         "07:00", "Prievidza",       "08:00", "Konzultácia vývoja softvéru",    127128, 116),
    @("2025-11-14", "Prievidza",      "16:30", "Vrbové",          "17:30", "Konzultácia vývoja softvéru",    127244, 116)
)

# Column A holds date-like text (e.g. "2025-11-03"); Excel's COM layer will
# silently re-interpret such strings as real dates unless the cell is first
# marked as Text. Mark the range as Text, write the values, then drop the
# style back to Normal so the saved cells stay plain/unstyled like the rest
# of the data cells.
$dateRange = $ws.Range("A2:A19")
$dateRange.NumberFormat = "@"

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

$dateRange.Style = "Normal"
